# "Removed Test Case Inter-Dependency"
# The product name / shortname on this test case were colliding with another
# test case's data, so both are made unique, and the Input sheet (rather than
# the Output sheet) is left as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Give this test case its own, non-colliding product name/short name.
$newProductName = "2570-MS-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-CASH-1st"
$newShortName   = "257d"

$wsInput.Range("B1").Value  = $newProductName
$wsOutput.Range("B1").Value = $newProductName

$wsInput.Range("B2").Value = $newShortName

# Leave the Input sheet active/selected (was previously left on the Output
# sheet), with the selection reset to the top cell.
$wsInput.Activate()
$wsInput.Range("B1").Select()
